$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values reflecting the new TPM-derived computation.
# Columns: E=Ligand-expressing cells, F=Ligand detection rate,
#          G=Ligand average expression value, H=Ligand total expression value,
#          M=Receptor average expression value, N=Receptor total expression value,
#          O=Receptor derived specificity (avg), P=Receptor derived specificity (total),
#          Q=Edge average expression weight, R=Edge total expression weight,
#          S=Edge average expression derived specificity, T=Edge total expression derived specificity

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.298672
$ws.Range("H2").Value = 3.896016
$ws.Range("M2").Value = 1.370876333333333
$ws.Range("N2").Value = 4.112629
$ws.Range("O2").Value = 0.01103063309339269
$ws.Range("P2").Value = 0.01103063309339269
$ws.Range("Q2").Value = 1.780318709562667
$ws.Range("R2").Value = 16.022868386064
$ws.Range("S2").Value = 0.01103063309339269
$ws.Range("T2").Value = 0.01103063309339269

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.298672
$ws.Range("H3").Value = 3.896016
$ws.Range("O3").Value = 0.7476219244149905
$ws.Range("P3").Value = 0.7476219244149904
$ws.Range("Q3").Value = 120.6644522074187
$ws.Range("R3").Value = 1085.980069866768
$ws.Range("S3").Value = 0.7476219244149905
$ws.Range("T3").Value = 0.7476219244149904

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.298672
$ws.Range("H4").Value = 3.896016
$ws.Range("M4").Value = 29.718484
$ws.Range("N4").Value = 89.155452
$ws.Range("O4").Value = 0.2391271080585153
$ws.Range("P4").Value = 0.2391271080585153
$ws.Range("Q4").Value = 38.594563053248
$ws.Range("R4").Value = 347.351067479232
$ws.Range("S4").Value = 0.2391271080585153
$ws.Range("T4").Value = 0.2391271080585153

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.298672
$ws.Range("H5").Value = 3.896016
$ws.Range("M5").Value = 0.275941
$ws.Range("N5").Value = 0.827823
$ws.Range("O5").Value = 0.002220334433101459
$ws.Range("P5").Value = 0.002220334433101458
$ws.Range("Q5").Value = 0.358356850352
$ws.Range("R5").Value = 3.225211653168
$ws.Range("S5").Value = 0.002220334433101459
$ws.Range("T5").Value = 0.002220334433101458
